$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 100002696
$ws.Range("J40").Value = 125002620
$ws.Range("L40").Value = 125002620
$ws.Range("N40").Value = -125002970
$ws.Range("H62").Value = 3930.6667
$ws.Range("I62").Value = 3930.6667
$ws.Range("K62").Value = 3930.6667
$ws.Range("M62").Value = -3306.6667
$ws.Range("H65").Value = 3930.6667
$ws.Range("I65").Value = 3930.6667
$ws.Range("K65").Value = 19653.3335
$ws.Range("M65").Value = -16533.3335
$ws.Range("H106").Value = 9026.833000000001
$ws.Range("I106").Value = 9026.833000000001
$ws.Range("K106").Value = 9026.833000000001
$ws.Range("M106").Value = -8395.833000000001
$ws.Range("H112").Value = 2196.3823
$ws.Range("J112").Value = 2232.6667
$ws.Range("L112").Value = 6698.000100000001
$ws.Range("N112").Value = -8914.000100000001
$ws.Range("H132").Value = 1994.322
$ws.Range("J132").Value = 2116.3333
$ws.Range("L132").Value = 6348.999899999999
$ws.Range("N132").Value = -11408.9999
$ws.Range("H138").Value = 10580.02
$ws.Range("I138").Value = 8056
$ws.Range("J138").Value = 10712.863
$ws.Range("K138").Value = 24168
$ws.Range("L138").Value = 32138.589
$ws.Range("M138").Value = -19028
$ws.Range("N138").Value = -42418.589

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1522.7778
$ws.Range("I2").Value = 1255.5454
$ws.Range("K2").Value = 1255.5454
$ws.Range("M2").Value = -1142.5454
$ws.Range("H32").Value = 18689.775
$ws.Range("I32").Value = 18689.775
$ws.Range("K32").Value = 18689.775
$ws.Range("M32").Value = -18402.775
$ws.Range("H116").Value = 1522.7778
$ws.Range("I116").Value = 1255.5454
$ws.Range("K116").Value = 1255.5454
$ws.Range("M116").Value = 1038.4546
$ws.Range("H122").Value = 6414.9546
$ws.Range("I122").Value = 6414.9546
$ws.Range("K122").Value = 19244.8638
$ws.Range("M122").Value = -16794.8638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1522.7778
$ws.Range("I3").Value = 1255.5454
$ws.Range("K3").Value = 1255.5454
$ws.Range("M3").Value = -1141.5454
$ws.Range("H82").Value = 35199.6
$ws.Range("I82").Value = 6499.5
$ws.Range("K82").Value = 6499.5
$ws.Range("M82").Value = -6116.5
$ws.Range("H85").Value = 35199.6
$ws.Range("I85").Value = 6499.5
$ws.Range("K85").Value = 6499.5
$ws.Range("M85").Value = -5173.5
$ws.Range("H99").Value = 1628.6316
$ws.Range("I99").Value = 1702.8667
$ws.Range("K99").Value = 1702.8667
$ws.Range("M99").Value = -204.8667
$ws.Range("H132").Value = 121208.43
$ws.Range("J132").Value = 121208.43
$ws.Range("L132").Value = 121208.43
$ws.Range("N132").Value = -131328.43
$ws.Range("H137").Value = 112608.164
$ws.Range("J137").Value = 116706.71
$ws.Range("L137").Value = 116706.71
$ws.Range("N137").Value = -126906.71
$ws.Range("H140").Value = 201998.8
$ws.Range("J140").Value = 239998.5
$ws.Range("L140").Value = 239998.5
$ws.Range("N140").Value = -250358.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83.2
$ws.Range("I7").Value = 99.666664
$ws.Range("J7").Value = 58.5
$ws.Range("K7").Value = 99.666664
$ws.Range("L7").Value = 58.5
$ws.Range("M7").Value = 13.333336
$ws.Range("N7").Value = -284.5
$ws.Range("H28").Value = 30699.5
$ws.Range("J28").Value = 30699.5
$ws.Range("L28").Value = 30699.5
$ws.Range("N28").Value = -31189.5
$ws.Range("H31").Value = 41671868
$ws.Range("I31").Value = 83337030
$ws.Range("K31").Value = 83337030
$ws.Range("M31").Value = -83336735
$ws.Range("H34").Value = 41671868
$ws.Range("I34").Value = 83337030
$ws.Range("K34").Value = 83337030
$ws.Range("M34").Value = -83336828
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 70000
$ws.Range("J80").Value = 70000
$ws.Range("L80").Value = 70000
$ws.Range("N80").Value = -72246
$ws.Range("H83").Value = 70000
$ws.Range("J83").Value = 70000
$ws.Range("L83").Value = 210000
$ws.Range("N83").Value = -221232
$ws.Range("H134").Value = 1674.6
$ws.Range("I134").Value = 1343.25
$ws.Range("K134").Value = 4029.75
$ws.Range("M134").Value = -1494.75
$ws.Range("H141").Value = 639895.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 639895.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 639895.75
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -650255.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 5572.8
$ws.Range("I17").Value = 3612.5
$ws.Range("K17").Value = 10837.5
$ws.Range("M17").Value = -10668.5
$ws.Range("H38").Value = 31.217392
$ws.Range("I38").Value = 36.555557
$ws.Range("J38").Value = 12
$ws.Range("K38").Value = 109.666671
$ws.Range("L38").Value = 36
$ws.Range("M38").Value = 237.333329
$ws.Range("N38").Value = -730

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H126").Value = 10870885
$ws.Range("I126").Value = 14821871
$ws.Range("K126").Value = 44465613
$ws.Range("M126").Value = -44463143

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9212.789000000001
$ws.Range("I40").Value = 5921.3184
$ws.Range("J40").Value = 13738.5625
$ws.Range("K40").Value = 5921.3184
$ws.Range("L40").Value = 13738.5625
$ws.Range("M40").Value = -5785.3184
$ws.Range("N40").Value = -14010.5625
$ws.Range("H46").Value = 1217
$ws.Range("J46").Value = 1136
$ws.Range("L46").Value = 1136
$ws.Range("N46").Value = -1512
$ws.Range("H122").Value = 3403.0208
$ws.Range("I122").Value = 3400.9575
$ws.Range("K122").Value = 10202.8725
$ws.Range("M122").Value = -7752.872499999999
$ws.Range("H124").Value = 99188.22
$ws.Range("J124").Value = 99188.22
$ws.Range("L124").Value = 99188.22
$ws.Range("N124").Value = -109008.22
$ws.Range("H132").Value = 5020.5713
$ws.Range("I132").Value = 3102.1428
$ws.Range("K132").Value = 9306.428400000001
$ws.Range("M132").Value = -6776.428400000001
$ws.Range("H136").Value = 10065.077
$ws.Range("I136").Value = 15414.667
$ws.Range("K136").Value = 46244.001
$ws.Range("M136").Value = -43694.001
$ws.Range("H137").Value = 114617.7
$ws.Range("J137").Value = 117233.4
$ws.Range("L137").Value = 117233.4
$ws.Range("N137").Value = -127433.4
